$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.630.79'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '3.445.47'
$ws.Range("E3").Value = '  -3.11%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '593.13'
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("D6").Value = '136.21'
$ws.Range("E6").Value = '  -7.30%  '
$ws.Range("D7").Value = '3.444.00'
$ws.Range("E7").Value = '  -3.11%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("D10").Value = '7.39'
$ws.Range("E10").Value = '  -5.78%  '
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").Value = '  -8.52%  '
$ws.Range("D12").Value = '0.379'
$ws.Range("E12").Value = '  -7.51%  '
$ws.Range("D13").Value = '4.025.47'
$ws.Range("E13").Value = '  -3.27%  '
$ws.Range("D14").Value = '0.0000182'
$ws.Range("E14").Value = '  -10.15%  '
$ws.Range("D15").Value = '26.68'
$ws.Range("E15").Value = '  -8.79%  '
$ws.Range("D16").Value = '3.465.58'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '65.583.96'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '0.115'
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("D19").Value = '9.90'
$ws.Range("E19").Value = '  -10.50%  '
$ws.Range("D20").Value = '5.84'
$ws.Range("E20").Value = '  -6.65%  '
$ws.Range("D21").Value = '13.76'
$ws.Range("E21").Value = '  -7.15%  '
$ws.Range("D22").Value = '393.32'
$ws.Range("E22").Value = '  -6.24%  '
$ws.Range("D23").Value = '0.552'
$ws.Range("E23").Value = '  -8.88%  '
$ws.Range("D24").Value = '73.51'
$ws.Range("E24").Value = '  -5.58%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = '3.587.03'
$ws.Range("E26").Value = '  -2.99%  '
$ws.Range("D27").Value = '0.0000106'
$ws.Range("E27").Value = '  -9.94%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '2.26'
$ws.Range("E29").Value = '  -8.84%  '
$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  -9.04%  '
$ws.Range("D31").Value = '8.20'
$ws.Range("E31").Value = '  -11.23%  '
$ws.Range("D32").Value = '3.450.74'
$ws.Range("E32").Value = '  -2.89%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '0.147'
$ws.Range("E34").Value = '  -5.99%  '
$ws.Range("D35").Value = '23.04'
$ws.Range("E35").Value = '  -6.64%  '
$ws.Range("D36").Value = '172.44'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").Value = '6.98'
$ws.Range("E37").Value = '  -8.89%  '
$ws.Range("D38").Value = '1.20'
$ws.Range("E38").Value = '  -9.77%  '
$ws.Range("D39").Value = '1.51'
$ws.Range("E39").Value = '  -6.66%  '
$ws.Range("D40").Value = '4.84'
$ws.Range("E40").Value = '  -9.72%  '
$ws.Range("D41").Value = '0.0770'
$ws.Range("E41").Value = '  -7.25%  '
$ws.Range("D42").Value = '0.825'
$ws.Range("E42").Value = '  -4.68%  '
$ws.Range("D43").Value = '43.66'
$ws.Range("E43").Value = '  -4.66%  '
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '4.42'
$ws.Range("E45").Value = '  -13.95%  '
$ws.Range("D46").Value = '1.62'
$ws.Range("E46").Value = '  -11.53%  '
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").Value = '22.49'
$ws.Range("E48").Value = '  -2.72%  '
$ws.Range("D49").Value = '6.57'
$ws.Range("E49").Value = '  -8.61%  '
$ws.Range("D50").Value = '2.09'
$ws.Range("E50").Value = '  -14.75%  '
$ws.Range("D51").Value = '2.200.76'
$ws.Range("E51").Value = '  -7.55%  '
